$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the latest GitHub Actions scrape run.
# Cells are formatted as Text first so values like "1.000" or "245.11" are not
# reinterpreted by Excel as numbers, matching the original inlineStr storage.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '28.957.98'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -1.49%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.834.61'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -1.84%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '245.11'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +0.52%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.6895'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -2.33%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.07708'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -2.83%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.3054'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -2.65%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '23.53'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -3.93%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07803'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.841.97'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -1.39%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.075'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -2.11%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '90.53'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -3.49%  '

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -3.05%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '6.441'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -1.33%  '

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -0.61%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '28.967.36'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -1.48%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '243.62'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -4.26%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '2.084.73'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -1.53%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '12.73'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -2.99%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.9996'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '7.485'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -2.25%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '163.55'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +1.59%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.1471'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -5.54%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.808'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -2.31%  '

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -3.36%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.554'
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +3.41%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.217'
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -2.54%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.158'
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -2.18%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.170'
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -3.46%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.05118'
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -3.38%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.7731'
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +3.12%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.850'
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -2.57%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.144'
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -2.74%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.682'
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -1.00%  '

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -1.79%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.241.43'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -3.78%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.697'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -2.40%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.9405'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +5.31%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '108.57'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -0.22%  '

$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = 'PaxDollar'
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '

$c = $ws.Range("B44")
$c.NumberFormat = "@"
$c.Value = 'FraxShare'
$c = $ws.Range("C44")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '5.738'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -4.60%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '9.612'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +0.19%  '

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -3.79%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.983.58'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -1.88%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.5172'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -0.12%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '64.45'
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -9.33%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.750'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -2.73%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.929'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -2.09%  '
